$wb = $excel.ActiveWorkbook

$about = $wb.Worksheets.Item("About")
$data  = $wb.Worksheets.Item("MSCdtRPbQL")

# ---------------------------------------------------------------------------
# Sheet "MSCdtRPbQL": clarify the variable name and let it wrap onto two
# lines (row made taller to fit the wrapped text).
# ---------------------------------------------------------------------------
$data.Range("A2").Value = "Change in Perc Share (dimensionless)"
$data.Range("A2").WrapText = $true
$data.Rows.Item(2).RowHeight = 30

# ---------------------------------------------------------------------------
# Sheet "About"
# ---------------------------------------------------------------------------

# The citation's URL moved; update the hyperlink target as well as the
# visible cell text, and keep the built-in "Hyperlink" cell style.
$citationCell = $about.Range("B6")
$newUrl = "https://ethz.ch/content/dam/ethz/special-interest/mtec/cepe/cepe-dam/documents/research/cepe-wp/CEPE_WP86.pdf"
$citationCell.Hyperlinks.Delete()
$about.Hyperlinks.Add($citationCell, $newUrl)
$citationCell.Value = $newUrl
$citationCell.Style = "Hyperlink"

# Insert a new row at row 10 (pushes the old rows 10.. down by one) and fill
# it with a new explanatory sentence about the "Change in Perc Share" metric.
$about.Rows.Item(10).Insert()
$about.Range("A10").Value = "This variable measures how a rebate program influences market shares of rebate-qualifying and non-qualifying components."
# The inserted row inherited the bold formatting of row 9 ("Note:") - the
# new note itself is plain (non-bold) body text.
$about.Range("A10").Font.Bold = $false
